$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Row 11 was a placeholder ("**") waiting for the next repair-log entry.
# Fill it in with the new entry: date, technician, defect reason, strike count.
# The source cells all store digit-looking values as plain TEXT (shared
# strings), so force Text format before writing, then drop the border that
# inheriting the column's style would otherwise add.

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "06/02/2018"
$ws.Range("A11").Borders.LineStyle = -4142

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "3012"
$ws.Range("B11").Borders.LineStyle = -4142

$ws.Range("C11").Value = "Зазубрини в місті відрізу контакту"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22222222222222"
$ws.Range("D11").Borders.LineStyle = -4142

# The "**" placeholder moves down to row 12, marking the next entry to fill.
$ws.Range("A12").Value = "**"
